$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Amelx"
$ws.Range("C2").Value = "Lamp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.069782333333333
$ws.Range("H2").Value = 3.209347
$ws.Range("I2").Value = 0.5403049320348338
$ws.Range("J2").Value = 0.5403049320348337
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 62.86930466666666
$ws.Range("N2").Value = 188.607914
$ws.Range("O2").Value = 0.110909371613604
$ws.Range("P2").Value = 0.110909371613604
$ws.Range("Q2").Value = 67.2564714413509
$ws.Range("R2").Value = 605.308242972158
$ws.Range("S2").Value = 0.05992488049171445
$ws.Range("T2").Value = 0.05992488049171445

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Amelx"
$ws.Range("C3").Value = "Lamp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.069782333333333
$ws.Range("H3").Value = 3.209347
$ws.Range("I3").Value = 0.5403049320348338
$ws.Range("J3").Value = 0.5403049320348337
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 241.00591
$ws.Range("N3").Value = 723.01773
$ws.Range("O3").Value = 0.42516477913962
$ws.Range("P3").Value = 0.42516477913962
$ws.Range("Q3").Value = 257.8238647469233
$ws.Range("R3").Value = 2320.41478272231
$ws.Range("S3").Value = 0.2297186270966375
$ws.Range("T3").Value = 0.2297186270966374

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Amelx"
$ws.Range("C4").Value = "Lamp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.069782333333333
$ws.Range("H4").Value = 3.209347
$ws.Range("I4").Value = 0.5403049320348338
$ws.Range("J4").Value = 0.5403049320348337
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 97.28925299999999
$ws.Range("N4").Value = 291.867759
$ws.Range("O4").Value = 0.1716304955525929
$ws.Range("P4").Value = 0.1716304955525929
$ws.Range("Q4").Value = 104.078324082597
$ws.Range("R4").Value = 936.704916743373
$ws.Range("S4").Value = 0.09273280323464855
$ws.Range("T4").Value = 0.09273280323464853

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Amelx"
$ws.Range("C5").Value = "Lamp2"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.069782333333333
$ws.Range("H5").Value = 3.209347
$ws.Range("I5").Value = 0.5403049320348338
$ws.Range("J5").Value = 0.5403049320348337
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 165.688484
$ws.Range("N5").Value = 497.065452
$ws.Range("O5").Value = 0.2922953536941831
$ws.Range("P5").Value = 0.2922953536941831
$ws.Range("Q5").Value = 177.2506130199827
$ws.Range("R5").Value = 1595.255517179844
$ws.Range("S5").Value = 0.1579286212118333
$ws.Range("T5").Value = 0.1579286212118333

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Amelx"
$ws.Range("C6").Value = "Lamp2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.663689
$ws.Range("H6").Value = 1.991067
$ws.Range("I6").Value = 0.3352031799963669
$ws.Range("J6").Value = 0.3352031799963669
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 62.86930466666666
$ws.Range("N6").Value = 188.607914
$ws.Range("O6").Value = 0.110909371613604
$ws.Range("P6").Value = 0.110909371613604
$ws.Range("Q6").Value = 41.72566594491533
$ws.Range("R6").Value = 375.530993504238
$ws.Range("S6").Value = 0.03717717405627886
$ws.Range("T6").Value = 0.03717717405627886

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Amelx"
$ws.Range("C7").Value = "Lamp2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.663689
$ws.Range("H7").Value = 1.991067
$ws.Range("I7").Value = 0.3352031799963669
$ws.Range("J7").Value = 0.3352031799963669
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 241.00591
$ws.Range("N7").Value = 723.01773
$ws.Range("O7").Value = 0.42516477913962
$ws.Range("P7").Value = 0.42516477913962
$ws.Range("Q7").Value = 159.95297140199
$ws.Range("R7").Value = 1439.57674261791
$ws.Range("S7").Value = 0.1425165859900536
$ws.Range("T7").Value = 0.1425165859900536

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Amelx"
$ws.Range("C8").Value = "Lamp2"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.663689
$ws.Range("H8").Value = 1.991067
$ws.Range("I8").Value = 0.3352031799963669
$ws.Range("J8").Value = 0.3352031799963669
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 97.28925299999999
$ws.Range("N8").Value = 291.867759
$ws.Range("O8").Value = 0.1716304955525929
$ws.Range("P8").Value = 0.1716304955525929
$ws.Range("Q8").Value = 64.56980703431699
$ws.Range("R8").Value = 581.128263308853
$ws.Range("S8").Value = 0.05753108789358145
$ws.Range("T8").Value = 0.05753108789358145

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Amelx"
$ws.Range("C9").Value = "Lamp2"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.663689
$ws.Range("H9").Value = 1.991067
$ws.Range("I9").Value = 0.3352031799963669
$ws.Range("J9").Value = 0.3352031799963669
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 165.688484
$ws.Range("N9").Value = 497.065452
$ws.Range("O9").Value = 0.2922953536941831
$ws.Range("P9").Value = 0.2922953536941831
$ws.Range("Q9").Value = 109.965624257476
$ws.Range("R9").Value = 989.6906183172839
$ws.Range("S9").Value = 0.09797833205645298
$ws.Range("T9").Value = 0.09797833205645298

$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Amelx"
$ws.Range("C10").Value = "Lamp2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.246489
$ws.Range("H10").Value = 0.7394670000000001
$ws.Range("I10").Value = 0.1244918879687994
$ws.Range("J10").Value = 0.1244918879687994
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 62.86930466666666
$ws.Range("N10").Value = 188.607914
$ws.Range("O10").Value = 0.110909371613604
$ws.Range("P10").Value = 0.110909371613604
$ws.Range("Q10").Value = 15.496592037982
$ws.Range("R10").Value = 139.469328341838
$ws.Range("S10").Value = 0.01380731706561074
$ws.Range("T10").Value = 0.01380731706561073

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Amelx"
$ws.Range("C11").Value = "Lamp2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.246489
$ws.Range("H11").Value = 0.7394670000000001
$ws.Range("I11").Value = 0.1244918879687994
$ws.Range("J11").Value = 0.1244918879687994
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 241.00591
$ws.Range("N11").Value = 723.01773
$ws.Range("O11").Value = 0.42516477913962
$ws.Range("P11").Value = 0.42516477913962
$ws.Range("Q11").Value = 59.40530574999001
$ws.Range("R11").Value = 534.6477517499101
$ws.Range("S11").Value = 0.05292956605292891
$ws.Range("T11").Value = 0.0529295660529289

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Amelx"
$ws.Range("C12").Value = "Lamp2"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.246489
$ws.Range("H12").Value = 0.7394670000000001
$ws.Range("I12").Value = 0.1244918879687994
$ws.Range("J12").Value = 0.1244918879687994
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 97.28925299999999
$ws.Range("N12").Value = 291.867759
$ws.Range("O12").Value = 0.1716304955525929
$ws.Range("P12").Value = 0.1716304955525929
$ws.Range("Q12").Value = 23.980730682717
$ws.Range("R12").Value = 215.826576144453
$ws.Range("S12").Value = 0.02136660442436292
$ws.Range("T12").Value = 0.02136660442436292

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Amelx"
$ws.Range("C13").Value = "Lamp2"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.246489
$ws.Range("H13").Value = 0.7394670000000001
$ws.Range("I13").Value = 0.1244918879687994
$ws.Range("J13").Value = 0.1244918879687994
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 165.688484
$ws.Range("N13").Value = 497.065452
$ws.Range("O13").Value = 0.2922953536941831
$ws.Range("P13").Value = 0.2922953536941831
$ws.Range("Q13").Value = 40.840388732676
$ws.Range("R13").Value = 367.563498594084
$ws.Range("S13").Value = 0.03638840042589683
$ws.Range("T13").Value = 0.03638840042589683
